$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - update "想去人数" (want-to-go count) values
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 293
$wsExpo.Range("F4").Value = 2431
$wsExpo.Range("F5").Value = 1785
$wsExpo.Range("F8").Value = 840

# Sheet "全部类型" (All types) - same underlying events, shifted by one row
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 293
$wsAll.Range("F4").Value = 2431
$wsAll.Range("F5").Value = 1785
$wsAll.Range("F9").Value = 840
